$d = $word.ActiveDocument

# The document contains 5 placeholder pictures (1x1 px) that need to become
# plain hyperlink runs whose visible text is the image's original URL.
# Shapes 1-2 -> Waterbodies images, shapes 3-5 -> Foreshore images.
$urls = @(
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C16_Waterbodies_1.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C17_Waterbodies_2.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C13_Foreshore_A.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C14_Foreshore_B.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C15_Foreshore_C.jpg?h=100%25&w=100%25"
)

# Always take the first remaining InlineShape: deleting it collapses the
# collection so the next picture becomes item 1, which keeps us walking the
# shapes in left-to-right / top-to-bottom document order.
$count = $d.InlineShapes.Count
for ($n = 0; $n -lt $count; $n++) {
    $shape = $d.InlineShapes.Item(1)
    $range = $shape.Range
    $shape.Delete()
    $url = $urls[$n]
    $d.Hyperlinks.Add($range, $url, $null, $null, $url) | Out-Null
}

Write-Output "done"
